# Powerpoint writer: Use table styles
#
# Apply the default table style (the one declared as `def` in the
# presentation's tableStyles part) to every table in the deck, so the
# generated tables match the reference document's color scheme instead
# of being left without an explicit <a:tableStyleId>.

$p = $ppt.ActivePresentation

# GUID of the default table style referenced by ppt/tableStyles.xml
# (def="{5C22544A-7EE6-4342-B048-85BDC9FD1C3A}").
$defaultTableStyleId = "{5C22544A-7EE6-4342-B048-85BDC9FD1C3A}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($defaultTableStyleId)
        }
    }
}
